$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for the Icam4-Itgb3 ligand-receptor pair sheet.
# Map of cell reference -> new numeric value, taken from the recomputed
# NATMI output using the new TPM-based expression matrix.
$updates = @{
    "G2" = 1.829945333333333
    "H2" = 5.489835999999999
    "I2" = 0.4190796720210465
    "J2" = 0.4190796720210465
    "M2" = 2.759544333333333
    "N2" = 8.278632999999999
    "O2" = 0.2574067337278401
    "P2" = 0.2574067337278401
    "Q2" = 5.049815274909776
    "R2" = 45.44833747418799
    "S2" = 0.107873929546672
    "T2" = 0.107873929546672
    "G3" = 1.829945333333333
    "H3" = 5.489835999999999
    "I3" = 0.4190796720210465
    "J3" = 0.4190796720210465
    "O3" = 0.6758254232987829
    "P3" = 0.6758254232987829
    "Q3" = 13.25836933759066
    "R3" = 119.325324038316
    "S3" = 0.2832246967395388
    "T3" = 0.2832246967395388
    "G4" = 1.829945333333333
    "H4" = 5.489835999999999
    "I4" = 0.4190796720210465
    "J4" = 0.4190796720210465
    "M4" = 0.5200313333333334
    "N4" = 1.560094
    "O4" = 0.0485078515798926
    "P4" = 0.0485078515798926
    "Q4" = 0.9516289116204445
    "R4" = 8.564660204583999
    "S4" = 0.02032865453054699
    "T4" = 0.02032865453054699
    "G5" = 1.829945333333333
    "H5" = 5.489835999999999
    "I5" = 0.4190796720210465
    "J5" = 0.4190796720210465
    "K5" = 3
    "L5" = 1
    "M5" = 0.1957573333333333
    "N5" = 0.587272
    "O5" = 0.01825999139348442
    "P5" = 0.01825999139348442
    "Q5" = 0.3582252185991111
    "R5" = 3.224026967392
    "S5" = 0.007652391204288582
    "T5" = 0.007652391204288582
    "I6" = 0.2833335737960661
    "J6" = 0.2833335737960661
    "M6" = 2.759544333333333
    "N6" = 8.278632999999999
    "O6" = 0.2574067337278401
    "P6" = 0.2574067337278401
    "Q6" = 3.414105489655666
    "R6" = 30.726949406901
    "S6" = 0.0729319697862813
    "T6" = 0.0729319697862813
    "I7" = 0.2833335737960661
    "J7" = 0.2833335737960661
    "O7" = 0.6758254232987829
    "P7" = 0.6758254232987829
    "S7" = 0.1914840324454833
    "T7" = 0.1914840324454833
    "I8" = 0.2833335737960661
    "J8" = 0.2833335737960661
    "M8" = 0.5200313333333334
    "N8" = 1.560094
    "O8" = 0.0485078515798926
    "P8" = 0.0485078515798926
    "Q8" = 0.6433822455686667
    "R8" = 5.790440210118
    "S8" = 0.01374390294530012
    "T8" = 0.01374390294530012
    "I9" = 0.2833335737960661
    "J9" = 0.2833335737960661
    "K9" = 3
    "L9" = 1
    "M9" = 0.1957573333333333
    "N9" = 0.587272
    "O9" = 0.01825999139348442
    "P9" = 0.01825999139348442
    "Q9" = 0.2421907770426667
    "R9" = 2.179716993384
    "S9" = 0.00517366861900135
    "T9" = 0.00517366861900135
    "E10" = 3
    "F10" = 1
    "G10" = 0.1530633333333333
    "H10" = 0.45919
    "I10" = 0.03505335944376924
    "J10" = 0.03505335944376924
    "M10" = 2.759544333333333
    "N10" = 8.278632999999999
    "O10" = 0.2574067337278401
    "P10" = 0.2574067337278401
    "Q10" = 0.422385054141111
    "R10" = 3.801465487269999
    "S10" = 0.009022970760608576
    "T10" = 0.009022970760608576
    "E11" = 3
    "F11" = 1
    "G11" = 0.1530633333333333
    "H11" = 0.45919
    "I11" = 0.03505335944376924
    "J11" = 0.03505335944376924
    "O11" = 0.6758254232987829
    "P11" = 0.6758254232987829
    "Q11" = 1.108978595376667
    "R11" = 9.980807358389999
    "S11" = 0.02368995148412973
    "T11" = 0.02368995148412973
    "E12" = 3
    "F12" = 1
    "G12" = 0.1530633333333333
    "H12" = 0.45919
    "I12" = 0.03505335944376924
    "J12" = 0.03505335944376924
    "M12" = 0.5200313333333334
    "N12" = 1.560094
    "O12" = 0.0485078515798926
    "P12" = 0.0485078515798926
    "Q12" = 0.07959772931777778
    "R12" = 0.71637956386
    "S12" = 0.001700363157274985
    "T12" = 0.001700363157274985
    "E13" = 3
    "F13" = 1
    "G13" = 0.1530633333333333
    "H13" = 0.45919
    "I13" = 0.03505335944376924
    "J13" = 0.03505335944376924
    "K13" = 3
    "L13" = 1
    "M13" = 0.1957573333333333
    "N13" = 0.587272
    "O13" = 0.01825999139348442
    "P13" = 0.01825999139348442
    "Q13" = 0.02996326996444445
    "R13" = 0.26966942968
    "S13" = 0.0006400740417559422
    "T13" = 0.0006400740417559422
    "G14" = 1.146373333333333
    "H14" = 3.43912
    "I14" = 0.2625333947391181
    "J14" = 0.2625333947391181
    "M14" = 2.759544333333333
    "N14" = 8.278632999999999
    "O14" = 0.2574067337278401
    "P14" = 0.2574067337278401
    "Q14" = 3.163468035884444
    "R14" = 28.47121232296
    "S14" = 0.0675778636342781
    "T14" = 0.0675778636342781
    "G15" = 1.146373333333333
    "H15" = 3.43912
    "I15" = 0.2625333947391181
    "J15" = 0.2625333947391181
    "O15" = 0.6758254232987829
    "P15" = 0.6758254232987829
    "Q15" = 8.305735026746666
    "R15" = 74.75161524072
    "S15" = 0.177426742629631
    "T15" = 0.177426742629631
    "G16" = 1.146373333333333
    "H16" = 3.43912
    "I16" = 0.2625333947391181
    "J16" = 0.2625333947391181
    "M16" = 0.5200313333333334
    "N16" = 1.560094
    "O16" = 0.0485078515798926
    "P16" = 0.0485078515798926
    "Q16" = 0.5961500530311111
    "R16" = 5.36535047728
    "S16" = 0.0127349309467705
    "T16" = 0.0127349309467705
    "G17" = 1.146373333333333
    "H17" = 3.43912
    "I17" = 0.2625333947391181
    "J17" = 0.2625333947391181
    "K17" = 3
    "L17" = 1
    "M17" = 0.1957573333333333
    "N17" = 0.587272
    "O17" = 0.01825999139348442
    "P17" = 0.01825999139348442
    "Q17" = 0.2244109867377778
    "R17" = 2.01969888064
    "S17" = 0.004793857528438546
    "T17" = 0.004793857528438546
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value2 = $updates[$ref]
}
